# "paises.xlsx" COVID-19 dashboard refresh: new pull at 20:58 replaces the
# 19:41 snapshot. Totals/new-cases/active/recovered/deaths-today/deaths move
# for the countries whose figures changed, and three country pairs swap rank
# order (Marruecos overtakes Bolivia, Libano overtakes Ghana, Santa Lucia
# overtakes Nueva Caledonia) because their "Casos totales" crossed — the two
# rows in each pair keep their stats, only the country label moves.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp in A1
$ws.Range("A1").Value = "Datos actualizados a 6 de Octubre de 2020 a las 20:58"

# Row 4
$ws.Cells.Item(4, 2).Value = 7695511
$ws.Cells.Item(4, 3).Value = 15867
$ws.Cells.Item(4, 4).Value = 4912156
$ws.Cells.Item(4, 5).Value = 2568032
$ws.Cells.Item(4, 7).Value = 291
$ws.Cells.Item(4, 8).Value = 215323

# Row 5
$ws.Cells.Item(5, 2).Value = 6741616
$ws.Cells.Item(5, 3).Value = 59543
$ws.Cells.Item(5, 4).Value = 5722286
$ws.Cells.Item(5, 5).Value = 915110
$ws.Cells.Item(5, 7).Value = 620
$ws.Cells.Item(5, 8).Value = 104220

# Row 14
$ws.Cells.Item(14, 2).Value = 634763
$ws.Cells.Item(14, 3).Value = 10489
$ws.Cells.Item(14, 4).Value = 99295
$ws.Cells.Item(14, 5).Value = 503103
$ws.Cells.Item(14, 7).Value = 66
$ws.Cells.Item(14, 8).Value = 32365

# Row 26
$ws.Cells.Item(26, 2).Value = 306651
$ws.Cells.Item(26, 3).Value = 1994
$ws.Cells.Item(26, 5).Value = 33321
$ws.Cells.Item(26, 7).Value = 14
$ws.Cells.Item(26, 8).Value = 9630

# Row 29
$ws.Cells.Item(29, 2).Value = 170872
$ws.Cells.Item(29, 3).Value = 1912
$ws.Cells.Item(29, 4).Value = 143733
$ws.Cells.Item(29, 5).Value = 17613
$ws.Cells.Item(29, 7).Value = 22
$ws.Cells.Item(29, 8).Value = 9526

# Row 33
$ws.Cells.Item(33, 1).Value = "Marruecos"
$ws.Cells.Item(33, 2).Value = 137248
$ws.Cells.Item(33, 3).Value = 2553
$ws.Cells.Item(33, 4).Value = 115354
$ws.Cells.Item(33, 5).Value = 19484
$ws.Cells.Item(33, 7).Value = 41
$ws.Cells.Item(33, 8).Value = 2410

# Row 34
$ws.Cells.Item(34, 1).Value = "Bolivia"
$ws.Cells.Item(34, 2).Value = 137107
$ws.Cells.Item(34, 3).Value = 239
$ws.Cells.Item(34, 4).Value = 98007
$ws.Cells.Item(34, 5).Value = 30971
$ws.Cells.Item(34, 7).Value = 28
$ws.Cells.Item(34, 8).Value = 8129

# Row 55
$ws.Cells.Item(55, 2).Value = 80003
$ws.Cells.Item(55, 3).Value = 566
$ws.Cells.Item(55, 4).Value = 34960
$ws.Cells.Item(55, 5).Value = 43805
$ws.Cells.Item(55, 7).Value = 8
$ws.Cells.Item(55, 8).Value = 1238

# Row 67
$ws.Cells.Item(67, 1).Value = "Libano"
$ws.Cells.Item(67, 2).Value = 46918
$ws.Cells.Item(67, 3).Value = 1261
$ws.Cells.Item(67, 4).Value = 20490
$ws.Cells.Item(67, 5).Value = 26004
$ws.Cells.Item(67, 7).Value = 10
$ws.Cells.Item(67, 8).Value = 424

# Row 68
$ws.Cells.Item(68, 1).Value = "Ghana"
$ws.Cells.Item(68, 2).Value = 46829
$ws.Cells.Item(68, 4).Value = 46060
$ws.Cells.Item(68, 5).Value = 466
$ws.Cells.Item(68, 8).Value = 303

# Row 70
$ws.Cells.Item(70, 2).Value = 42432
$ws.Cells.Item(70, 3).Value = 475
$ws.Cells.Item(70, 4).Value = 35599
$ws.Cells.Item(70, 5).Value = 6484
$ws.Cells.Item(70, 7).Value = 10
$ws.Cells.Item(70, 8).Value = 349

# Row 105
$ws.Cells.Item(105, 2).Value = 10621
$ws.Cells.Item(105, 3).Value = 54
$ws.Cells.Item(105, 4).Value = 9466
$ws.Cells.Item(105, 5).Value = 1121

# Row 112
$ws.Cells.Item(112, 2).Value = 8838
$ws.Cells.Item(112, 3).Value = 11
$ws.Cells.Item(112, 4).Value = 7013
$ws.Cells.Item(112, 5).Value = 1596

# Row 138
$ws.Cells.Item(138, 2).Value = 4108
$ws.Cells.Item(138, 3).Value = 14
$ws.Cells.Item(138, 4).Value = 3643
$ws.Cells.Item(138, 5).Value = 434

# Row 207
$ws.Cells.Item(207, 1).Value = "Santa Lucia"

# Row 208
$ws.Cells.Item(208, 1).Value = "Nueva Caledonia"
